# COVID-19 Bangladesh DataSheet update
# Adds one new day (2020-04-19, Excel serial 43940) of data to each of the
# three sheets (Confirmed, Recoverd, Death) and updates the active
# sheet/selection state to match the author's final view.

$wb = $excel.ActiveWorkbook

function Add-DayRow {
    param($ws, $newRow, $dateSerial, $formula, $newCaseValue)

    $srcRow = $newRow - 1

    # Copy formatting (styles) down from the last existing row so the new
    # row visually matches the rest of the table.
    $ws.Range("A" + $srcRow + ":C" + $srcRow).Copy() | Out-Null
    $ws.Range("A" + $newRow + ":C" + $newRow).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($newRow, 1).Value = $dateSerial
    $ws.Cells.Item($newRow, 2).Formula = $formula
    $ws.Cells.Item($newRow, 3).Value = $newCaseValue
}

# --- Confirmed (sheet1) ---
$wsConfirmed = $wb.Worksheets.Item("Confirmed")
Add-DayRow $wsConfirmed 44 43940 "=SUM(B43+C44)" 312

# --- Recoverd (sheet2) ---
$wsRecoverd = $wb.Worksheets.Item("Recoverd")
Add-DayRow $wsRecoverd 44 43940 "=SUM(B43+C44)" 9

# --- Death (sheet3) ---
$wsDeath = $wb.Worksheets.Item("Death")
Add-DayRow $wsDeath 44 43940 "=SUM(B43+C44)" 7

# --- Update selections / active sheet to match the saved workbook state ---
# Death sheet loses focus, selection left at B48
$wsDeath.Range("B48").Select() | Out-Null

# Recoverd sheet selection moves to C44
$wsRecoverd.Range("C44").Select() | Out-Null

# Confirmed sheet becomes the active tab, selection at C45
$wsConfirmed.Range("C45").Select() | Out-Null
